$d = $word.ActiveDocument

# Remove the telephone number mention from the Domain Model description
# (as well as the now-superfluous wording around it):
#   "...cognome, telefono, data di nascita (opzionale) ed un'immagine di profilo (opzionale)."
# becomes
#   "...cognome, data di nascita ed un'immagine di profilo (opzionale)."

$found = $d.Content.Find.Execute(
    "cognome, telefono, data di nascita (opzionale) ed", $true, $false, $false, $false, $false,
    $true, 1, $false, "cognome, data di nascita ed", 2)

if (-not $found) {
    throw "Expected text not found - telephone sentence could not be updated."
}
